# Rename the respondents in the availability sheet (column A, rows 2-14)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Eli"
$ws.Range("A3").Value = "Aaron"
$ws.Range("A4").Value = "Adrian"
$ws.Range("A5").Value = "Saeda"
$ws.Range("A6").Value = "Modupe"
$ws.Range("A7").Value = "Gabby"
$ws.Range("A8").Value = "Ngozi"
$ws.Range("A9").Value = "Ayomide"
$ws.Range("A10").Value = "Michael"
$ws.Range("A11").Value = "Leena "
$ws.Range("A12").Value = "Allen"
$ws.Range("A13").Value = "Morgan "
$ws.Range("A14").Value = "Leonard"

# Restore the view: scrolled so row 7 is at the top, with B13 selected
$ws.Range("B13").Select()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
